$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peru Liga 1")

# Row 61
$ws.Range("B61").Value = 6905578
$ws.Range("C61").Value = "Peru Liga 1"
$ws.Range("D61").Value = 45130.72916666666
$ws.Range("E61").Value = "AD Tarma"
$ws.Range("F61").Value = "Atletico Grau"
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = "H"
$ws.Range("L61").Value = 1.75
$ws.Range("M61").Value = 3.6
$ws.Range("N61").Value = 4
$ws.Range("O61").Value = 1.571
$ws.Range("P61").Value = 4.2
$ws.Range("Q61").Value = 5.75
$ws.Range("R61").Value = -1
$ws.Range("S61").Value = 1.975
$ws.Range("T61").Value = 1.825
$ws.Range("U61").Value = 2.5
$ws.Range("V61").Value = 1.8
$ws.Range("W61").Value = 2
$ws.Range("X61").Value = 0.571
$ws.Range("Y61").Value = -1
$ws.Range("Z61").Value = -1
$ws.Range("AA61").Value = 0
$ws.Range("AB61").Value = 0
$ws.Range("AC61").Value = -1
$ws.Range("AD61").Value = 1

# Row 62
$ws.Range("B62").Value = 6905571
$ws.Range("C62").Value = "Peru Liga 1"
$ws.Range("D62").Value = 45130.72916666666
$ws.Range("E62").Value = "FBC Melgar"
$ws.Range("F62").Value = "Sporting Cristal"
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 1
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = "D"
$ws.Range("L62").Value = 2.1
$ws.Range("M62").Value = 3.4
$ws.Range("N62").Value = 3
$ws.Range("O62").Value = 1.75
$ws.Range("P62").Value = 3.8
$ws.Range("Q62").Value = 4.75
$ws.Range("R62").Value = -0.75
$ws.Range("S62").Value = 1.95
$ws.Range("T62").Value = 1.85
$ws.Range("U62").Value = 2.5
$ws.Range("V62").Value = 1.95
$ws.Range("W62").Value = 1.85
$ws.Range("X62").Value = -1
$ws.Range("Y62").Value = 2.8
$ws.Range("Z62").Value = -1
$ws.Range("AA62").Value = -1
$ws.Range("AB62").Value = 0.8500000000000001
$ws.Range("AC62").Value = -1
$ws.Range("AD62").Value = 0.8500000000000001

# Row 156
$ws.Range("B156").Value = 7211640
$ws.Range("C156").Value = "Peru Liga 1"
$ws.Range("D156").Value = 45198.70833333334
$ws.Range("E156").Value = "UTC Cajamarca"
$ws.Range("F156").Value = "Sport Boys"
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 1
$ws.Range("I156").Value = 0
$ws.Range("J156").Value = 0
$ws.Range("K156").Value = "D"
$ws.Range("L156").Value = 1.615
$ws.Range("M156").Value = 3.75
$ws.Range("N156").Value = 5
$ws.Range("O156").Value = 1.5
$ws.Range("P156").Value = 4.2
$ws.Range("Q156").Value = 6.5
$ws.Range("R156").Value = -1
$ws.Range("S156").Value = 1.8
$ws.Range("T156").Value = 2.05
$ws.Range("U156").Value = 2.5
$ws.Range("V156").Value = 1.875
$ws.Range("W156").Value = 1.975
$ws.Range("X156").Value = -1
$ws.Range("Y156").Value = 3.2
$ws.Range("Z156").Value = -1
$ws.Range("AA156").Value = -1
$ws.Range("AB156").Value = 1.05
$ws.Range("AC156").Value = -1
$ws.Range("AD156").Value = 0.9750000000000001

# Row 157
$ws.Range("B157").Value = 7211641
$ws.Range("C157").Value = "Peru Liga 1"
$ws.Range("D157").Value = 45198.70833333334
$ws.Range("E157").Value = "Sport Huancayo"
$ws.Range("F157").Value = "Deportivo Municipal"
$ws.Range("G157").Value = 2
$ws.Range("H157").Value = 0
$ws.Range("I157").Value = 0
$ws.Range("J157").Value = 0
$ws.Range("K157").Value = "H"
$ws.Range("L157").Value = 1.125
$ws.Range("M157").Value = 7
$ws.Range("N157").Value = 17
$ws.Range("O157").Value = 1.166
$ws.Range("P157").Value = 6.5
$ws.Range("Q157").Value = 12
$ws.Range("R157").Value = -2
$ws.Range("S157").Value = 1.775
$ws.Range("T157").Value = 2.025
$ws.Range("U157").Value = 3.5
$ws.Range("V157").Value = 1.9
$ws.Range("W157").Value = 1.9
$ws.Range("X157").Value = 0.1659999999999999
$ws.Range("Y157").Value = -1
$ws.Range("Z157").Value = -1
$ws.Range("AA157").Value = 0
$ws.Range("AB157").Value = 0
$ws.Range("AC157").Value = -1
$ws.Range("AD157").Value = 0.8999999999999999

# Row 187
$ws.Range("B187").Value = 7384625
$ws.Range("C187").Value = "Peru Liga 1"
$ws.Range("D187").Value = 45228.70833333334
$ws.Range("E187").Value = "AD Tarma"
$ws.Range("F187").Value = "Carlos Manucci"
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0
$ws.Range("I187").Value = 0
$ws.Range("J187").Value = 0
$ws.Range("K187").Value = "D"
$ws.Range("L187").Value = 1.5
$ws.Range("M187").Value = 3.75
$ws.Range("N187").Value = 7
$ws.Range("O187").Value = 1.363
$ws.Range("P187").Value = 4.333
$ws.Range("Q187").Value = 9.5
$ws.Range("R187").Value = -1.25
$ws.Range("S187").Value = 1.875
$ws.Range("T187").Value = 1.925
$ws.Range("U187").Value = 2.5
$ws.Range("V187").Value = 1.8
$ws.Range("W187").Value = 2
$ws.Range("X187").Value = -1
$ws.Range("Y187").Value = 3.333
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = -1
$ws.Range("AB187").Value = 0.925
$ws.Range("AC187").Value = -1
$ws.Range("AD187").Value = 1

# Row 188
$ws.Range("B188").Value = 7384628
$ws.Range("C188").Value = "Peru Liga 1"
$ws.Range("D188").Value = 45228.70833333334
$ws.Range("E188").Value = "Deportivo Binacional"
$ws.Range("F188").Value = "FBC Melgar"
$ws.Range("G188").Value = 1
$ws.Range("H188").Value = 2
$ws.Range("I188").Value = 1
$ws.Range("J188").Value = 1
$ws.Range("K188").Value = "A"
$ws.Range("L188").Value = 2.75
$ws.Range("M188").Value = 3.3
$ws.Range("N188").Value = 2.375
$ws.Range("O188").Value = 3.3
$ws.Range("P188").Value = 3.6
$ws.Range("Q188").Value = 2
$ws.Range("R188").Value = 0.5
$ws.Range("S188").Value = 1.8
$ws.Range("T188").Value = 2
$ws.Range("U188").Value = 2.75
$ws.Range("V188").Value = 1.975
$ws.Range("W188").Value = 1.875
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = -1
$ws.Range("Z188").Value = 1
$ws.Range("AA188").Value = -1
$ws.Range("AB188").Value = 1
$ws.Range("AC188").Value = 0.4875
$ws.Range("AD188").Value = -0.5
